$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "vendor"
$ws.Range("B1").Value = "invoiceNumber"
$ws.Range("A2").Value = "Bank Account: 7387324"
$ws.Range("B2").Value = "INVOICE # 850888"
$ws.Range("A3").Value = "Bank Account: 7387324"
$ws.Range("B3").Value = "INVOICE # 850888"
